$p = $ppt.ActivePresentation

# --- Slide 3: fix "Pour Administrateurs" -> "Pour Administrateurs & Agents" ---
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(6).TextFrame.TextRange.Text = "Pour Administrateurs & Agents"

# --- Slide 5: remove the standalone technology-name text boxes ---
# (Flutter, Alpine.js + Tailwind, Django + DRF, PostgreSQL)
$s5 = $p.Slides.Item(5)

$namesToRemove = @("Flutter", "Alpine.js + Tailwind", "Django + DRF", "PostgreSQL")

for ($i = $s5.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $t = $shp.TextFrame.TextRange.Text
        if ($namesToRemove -contains $t) {
            $shp.Delete()
        }
    }
}
